# Thesis roadmap figure: thicken several thin connector/rectangle outlines
# and switch the "sysDot" dashed connectors to a regular "dash" style.
#
# Mapping (OOXML shape id -> Shapes collection index on slide 1):
#   Group A - line weight 0.25pt (w="3175") -> 0.5pt (w="6350"), dash untouched
#     24,25,26,27,28,29       -> 10-15  (Rectangle 23-28)
#     42,43,44,45,46,47,48,49 -> 28-35  (Straight Arrow Connector 41-48)
#     50,51,52,53             -> 36-39  (Straight Connector 49-52)
#     54                      -> 40     (Straight Arrow Connector 53)
#
#   Group B - line weight 1pt (w="12700") -> 0.75pt (w="9525"), dash sysDot -> dash
#             (no head/tail-end markers, so element order is unaffected)
#     39,40        -> 25,26  (Straight Connector 38,39)
#     453,454      -> 55,56  (Straight Connector 452,453)
#     456,457      -> 58,59  (Straight Connector 455,456)
#
#   Group C - same as Group B but the shape also carries head/tail-end
#             markers; re-touching the arrowhead properties after the dash
#             style keeps <a:prstDash> emitted ahead of <a:headEnd>/<a:tailEnd>
#             (matching the canonical element order) instead of after them.
#     41   -> 27  (Straight Connector 40)
#     455  -> 57  (Straight Connector 454)
#     458  -> 60  (Straight Connector 457)
#
#   Group D - dash sysDot -> dash only (weight stays 1.5pt / w="19050"),
#             also carries head/tail-end markers needing the same fix-up.
#     58   -> 44  (Straight Connector 57)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Group A
$thickenIdx = @(10,11,12,13,14,15,28,29,30,31,32,33,34,35,36,37,38,39,40)
foreach ($i in $thickenIdx) {
    $sh = $s.Shapes.Item($i)
    $sh.Line.Weight = 0.5
}

# Group B
$thinDashIdx = @(25,26,55,56,58,59)
foreach ($i in $thinDashIdx) {
    $sh = $s.Shapes.Item($i)
    $sh.Line.Weight = 0.75
    $sh.Line.DashStyle = 4
}

# Group C
$thinDashArrowIdx = @(27,57,60)
foreach ($i in $thinDashArrowIdx) {
    $sh = $s.Shapes.Item($i)
    $sh.Line.Weight = 0.75
    $sh.Line.DashStyle = 4
    $sh.Line.BeginArrowheadStyle = $sh.Line.BeginArrowheadStyle
    $sh.Line.EndArrowheadStyle = $sh.Line.EndArrowheadStyle
}

# Group D
$dashOnlyArrowIdx = @(44)
foreach ($i in $dashOnlyArrowIdx) {
    $sh = $s.Shapes.Item($i)
    $sh.Line.DashStyle = 4
    $sh.Line.BeginArrowheadStyle = $sh.Line.BeginArrowheadStyle
    $sh.Line.EndArrowheadStyle = $sh.Line.EndArrowheadStyle
}
